$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 13212.71166189251
$ws.Range("D5").Value = 13212.71166189251

$ws.Range("D9").Value = 7094.864641774657
$ws.Range("D10").Value = 7094.864641774657

$ws.Range("D14").Value = 7004.888338107765
$ws.Range("D15").Value = 7004.888338107765
